# Automatische test-sync: 2025-06-23 18:28:50
# Adds a new log entry (row 13) to the "Logs" sheet and bumps the
# "IT / Technisch probleem" tally on the "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A13").Value = "Kan mijn wachtwoord niet resetten"
$ws.Range("B13").Value = "mailmind.test@zohomail.eu"
$ws.Range("C13").Value = "Ik krijg geen e-mail bij wachtwoord resetten."
$ws.Range("D13").Value = "IT / Technisch probleem"
$ws.Range("E13").Value = "Beste klant,`nBedankt voor het contact opnemen. Het spijt me te horen dat u geen e-mail heeft ontvangen bij het resetten van uw wachtwoord. Om dit probleem op te lossen, hebben we wat meer informatie nodig. Kunt u ons uw gebruikersnaam en het e-mailadres dat is gekoppeld aan uw account doorgeven? Op die manier kunnen we verder onderzoeken wat er aan de hand is en u helpen uw wachtwoord te resetten.`nWe kijken uit naar uw antwoord.`nMet vriendelijke groet,`n[Naam] `nE-mailassistent bij [Bedrijfsnaam]"
$ws.Range("F13").Value = "2025-06-23 18:28:44"
$ws.Range("G13").Value = "Ja"
$ws.Rows.Item(13).AutoFit()

$dFc = $ws.Range("D2:D12").FormatConditions.Item(1)
$dFc.ModifyAppliesToRange($ws.Range("D2:D13"))

$gFc = $ws.Range("G2:G12").FormatConditions.Item(1)
$gFc.ModifyAppliesToRange($ws.Range("G2:G13"))

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 4
